# Apply the "Dayton" trace-matrix update:
#  - refresh a few Lines/Path references in the existing rows
#  - add the new trace-matrix row for the pseudocode doc (row 17)
#  - add the new trace-matrix row for the WIP networking note (row 18)
#  - move the selection/cursor the way the author left it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated Lines / Path values in the existing rows (8-16) ---
$ws.Range("G11").Value = "16-100"
$ws.Range("G12").Value = "116-191"
$ws.Range("G13").Value = "240-270"
$ws.Range("D14").Value = "../game/Blueprints/Test_Rooms"

# --- new row 17: pseudocode documentation entry ---
$ws.Range("D17").Value = "Documentation/Pseudocode"
$ws.Range("A17").Value = "Pseudocode for above algorithms"
$ws.Range("E17").Value = "Dayton - Pseudocode.docx"
$ws.Range("B17").Value = "Dayton"
$ws.Range("C17").Value = "Design"
$ws.Range("G17").Value = "Entire File"
$ws.Rows.Item(17).RowHeight = 43.2

# --- new row 18: work-in-progress networking note ---
$ws.Range("C18").Value = "Code - WIP"
$ws.Range("F18").Value = "RoomActorBase, LootingLootsGameModeBase, DoorActor, LootActor, AssetTemplate"
$ws.Range("G18").Value = "Yes"
$ws.Range("A18").Value = "Attempted networking of level layout"
$ws.Range("B18").Value = "Dayton"
$ws.Range("D18").Value = "../game/src"
$ws.Rows.Item(18).RowHeight = 72

# --- leave the selection where the author left it, scrolled back to the top ---
[void]$ws.Range("A1").Select()
[void]$ws.Range("C25").Select()
